# Append five new product rows (18-22) to Sheet1, continuing the existing
# product listing table (rows 2-17). Each row has:
#   A = product code, B = product name, C = "IAP01N", D = "1",
#   E = running sequence number, F = "RT"
# All of these must be stored as shared strings (t="s") using the same
# bordered cell style as the existing rows, exactly like the pre-existing
# data, even though most of the column A/E values look numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @("20140983", "HS MILD ORGNL 16'S",   "17"),
    @("20140982", "HS MILD MANGO ICE 16", "18"),
    @("20140987", "DJAVA FINE CUT 16'S",  "19"),
    @("20140988", "KING DJAVA FLTR 20'S", "20"),
    @("20140984", "SEN FILTER 12'S",      "21")
)

$startRow = 18
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $prevRow = $r - 1

    # Copy the row above (formatting, borders, style) into the new row so
    # it matches the existing table's appearance (style index unchanged).
    $ws.Range("A" + $prevRow + ":F" + $prevRow).Copy() | Out-Null
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    # Enter each value as a text-producing formula ( ="..." ) so digit-only
    # strings (product codes, sequence numbers) are not auto-coerced into
    # numeric cells, then immediately flatten the formulas down to plain
    # values (xlPasteValues) so the saved cells are ordinary shared-string
    # values with no formula left behind - matching the target file, which
    # stores every column as a shared string.
    $ws.Cells.Item($r, 1).Formula = "=""" + $data[0] + """"
    $ws.Cells.Item($r, 2).Value   = $data[1]
    $ws.Cells.Item($r, 3).Value   = "IAP01N"
    $ws.Cells.Item($r, 4).Formula = "=""1"""
    $ws.Cells.Item($r, 5).Formula = "=""" + $data[2] + """"
    $ws.Cells.Item($r, 6).Value   = "RT"

    $rowRange = $ws.Range("A" + $r + ":F" + $r)
    $rowRange.Copy() | Out-Null
    $rowRange.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = 0
}
